# Insert a new data row at row 271 (pushes existing rows 271-324 down to 272-325)
# and populate it with the new "Choclo" price record for
# "Dulce o Americano" / Región de Arica y Parinacota.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(271).Insert()

$ws.Cells.Item(271, 1).Value = 9
$ws.Cells.Item(271, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(271, 3).Value = "Metropolitana"
$ws.Cells.Item(271, 4).Value = 44522
$ws.Cells.Item(271, 5).Value = 13
$ws.Cells.Item(271, 6).Value = 100112024
$ws.Cells.Item(271, 7).Value = "Choclo"
$ws.Cells.Item(271, 8).Value = "Dulce o Americano"
$ws.Cells.Item(271, 9).Value = "Primera"
$ws.Cells.Item(271, 10).Value = 52
$ws.Cells.Item(271, 11).Value = 22000
$ws.Cells.Item(271, 12).Value = 25000
$ws.Cells.Item(271, 13).Value = 23500
$ws.Cells.Item(271, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(271, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(271, 16).Value = 336
$ws.Cells.Item(271, 17).Value = 70
$ws.Cells.Item(271, 18).Value = "Hortaliza"
